# Add data from table 4 (pT-stage factor rows) to the
# "ldsurvival-inputs" sheet, rows 24-30, columns E/F/H/I/L,
# and freeze the header row + first column (commit: "Add data from
# table 4 to excel file.")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ldsurvival-inputs")
$ws.Activate()

# --- Row 24 : pT1a (baseline level) -----------------------------------
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = ":v-radio"
$ws.Range("H24").Value = "[’p ‘TBD:]"
$ws.Range("I24").Value = 50
$ws.Range("L24").Value = "pT1a"

# --- Row 25 : pT1b ------------------------------------------------------
$ws.Range("E25").Value = "1,24"
$ws.Range("I25").Value = 50
$ws.Range("L25").Value = "pT1b"

# --- Row 26 : pT2 --------------------------------------------------------
$ws.Range("E26").Value = "1,64"
$ws.Range("I26").Value = 50
$ws.Range("L26").Value = "pT2"

# --- Row 27 : pT3a -------------------------------------------------------
$ws.Range("E27").Value = "1,8"
$ws.Range("I27").Value = 50
$ws.Range("L27").Value = "pT3a"

# --- Row 28 : pT3b -------------------------------------------------------
$ws.Range("E28").Value = "2,01"
$ws.Range("I28").Value = 50
$ws.Range("L28").Value = "pT3b"

# --- Row 29 : pT3c -------------------------------------------------------
$ws.Range("E29").Value = "2,01"
$ws.Range("I29").Value = 50
$ws.Range("L29").Value = "pT3c"

# --- Row 30 : pT4 --------------------------------------------------------
$ws.Range("E30").Value = "2,01"
$ws.Range("I30").Value = -50
$ws.Range("L30").Value = "pT4"

# --- View: freeze header row (1) + first column (A), then leave the
#     final selection on L30, matching the post-edit sheetView ---------
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("L30").Select()
